# Auto-generated script to apply the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings are
# stored verbatim (matching the source data which uses localized grouping
# like "26.765.27") instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.765.27"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "1.838.15"
$ws.Range("E3").Value = "  +1.78%  "

$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "309.46"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").Value = "0.4677"
$ws.Range("E7").Value = "  +3.38%  "

$ws.Range("D8").Value = "0.3617"
$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("D9").Value = "0.07161"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "0.9360"
$ws.Range("E10").Value = "  +5.11%  "

$ws.Range("D11").Value = "19.53"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "0.07675"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").Value = "1.853.57"
$ws.Range("E13").Value = "  +2.57%  "

$ws.Range("D14").Value = "5.275"
$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").Value = "6.367"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").Value = "87.93"
$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").Value = "0.000008555"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "26.773.56"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("D21").Value = "14.31"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "5.025"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("D24").Value = "1.918"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D25").Value = "152.20"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "17.99"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("D27").Value = "2.001"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").Value = "113.78"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").Value = "4.898"
$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").Value = "0.08838"
$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("D31").Value = "3.159"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").Value = "2.852"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "1.180"
$ws.Range("E33").Value = "  +6.70%  "

$ws.Range("D34").Value = "0.7423"
$ws.Range("E34").Value = "  +3.09%  "

$ws.Range("D35").Value = "4.448"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").Value = "1.086"
$ws.Range("E36").Value = "  +1.30%  "

$ws.Range("D37").Value = "2.981"
$ws.Range("E37").Value = "  +2.78%  "

$ws.Range("D38").Value = "0.01928"
$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("D39").Value = "0.05147"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").Value = "6.911"
$ws.Range("E40").Value = "  +1.97%  "

$ws.Range("D41").Value = "0.5094"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").Value = "0.1508"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("D43").Value = "8.134"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("D44").Value = "0.4684"
$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("D47").Value = "99.77"
$ws.Range("E47").Value = "  -0.59%  "

$ws.Range("D48").Value = "1.580"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").Value = "64.04"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").Value = "36.04"
$ws.Range("E51").Value = "  -0.23%  "

# Row 45: coin identity swap -> EnergySwap now ranks here
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.26"
$ws.Range("E45").Value = "  +2.80%  "

# Row 46: coin identity swap -> PaxDollar now ranks here
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  +0.49%  "

# Restore default styling on column D now that the text values are
# committed, so cells don't retain an explicit text number-format style.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Applied cryptos update"